$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (values changed in place) ---
# Row 2
$ws.Range("D2").Value = 44187
$ws.Range("M2").Value = 80
$ws.Range("N2").Value = 2800
$ws.Range("P2").Value = 2900
$ws.Range("R2").Value = "Provincia de Linares"
$ws.Range("S2").Value = 1450

# Row 3
$ws.Range("D3").Value = 44187
$ws.Range("M3").Value = 65
$ws.Range("N3").Value = 1400
$ws.Range("O3").Value = 1500
$ws.Range("P3").Value = 1446
$ws.Range("Q3").Value = "`$/envase 1 kilo"
$ws.Range("S3").Value = 1446
$ws.Range("T3").Value = 1

# Row 4
$ws.Range("D4").Value = 44956
$ws.Range("L4").Value = "Primera"
$ws.Range("N4").Value = 3000
$ws.Range("O4").Value = 3000
$ws.Range("P4").Value = 3000
$ws.Range("S4").Value = 1500

# Row 5
$ws.Range("D5").Value = 44971
$ws.Range("M5").Value = 30
$ws.Range("N5").Value = 3000
$ws.Range("P5").Value = 3000
$ws.Range("S5").Value = 1500

# Row 6
$ws.Range("D6").Value = 44974
$ws.Range("M6").Value = 40

# Row 7
$ws.Range("D7").Value = 44974

# Row 8
$ws.Range("D8").Value = 44935
$ws.Range("M8").Value = 50

# Row 9
$ws.Range("D9").Value = 44951
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 60
$ws.Range("N9").Value = 2800
$ws.Range("O9").Value = 3000
$ws.Range("P9").Value = 2900
$ws.Range("S9").Value = 1450

# Row 10
$ws.Range("D10").Value = 44932
$ws.Range("M10").Value = 60
$ws.Range("N10").Value = 3000
$ws.Range("P10").Value = 3000
$ws.Range("R10").Value = "Provincia de Diguillín"
$ws.Range("S10").Value = 1500

# Row 11
$ws.Range("D11").Value = 44949
$ws.Range("M11").Value = 60
$ws.Range("N11").Value = 2800
$ws.Range("O11").Value = 3000
$ws.Range("P11").Value = 2900
$ws.Range("Q11").Value = "`$/bandeja 2 kilos"
$ws.Range("S11").Value = 1450
$ws.Range("T11").Value = 2

# Row 12
$ws.Range("D12").Value = 44594
$ws.Range("M12").Value = 120
$ws.Range("O12").Value = 2800
$ws.Range("P12").Value = 2650
$ws.Range("R12").Value = "Provincia de Linares"
$ws.Range("S12").Value = 1325

# Row 13
$ws.Range("D13").Value = 44967
$ws.Range("M13").Value = 50
$ws.Range("N13").Value = 3000
$ws.Range("P13").Value = 3000
$ws.Range("S13").Value = 1500

# Row 14
$ws.Range("D14").Value = 44967
$ws.Range("L14").Value = "Segunda"
$ws.Range("M14").Value = 30
$ws.Range("N14").Value = 2500
$ws.Range("O14").Value = 2500
$ws.Range("P14").Value = 2500
$ws.Range("R14").Value = "Provincia de Diguillín"
$ws.Range("S14").Value = 1250

# Row 15
$ws.Range("D15").Value = 44181
$ws.Range("M15").Value = 65
$ws.Range("N15").Value = 3600
$ws.Range("O15").Value = 3800
$ws.Range("P15").Value = 3692
$ws.Range("S15").Value = 1846

# Row 16
$ws.Range("D16").Value = 44181
$ws.Range("M16").Value = 80
$ws.Range("N16").Value = 1800
$ws.Range("O16").Value = 2000
$ws.Range("P16").Value = 1875
$ws.Range("Q16").Value = "`$/envase 1 kilo"
$ws.Range("S16").Value = 1875
$ws.Range("T16").Value = 1

# Row 17
$ws.Range("D17").Value = 44965
$ws.Range("M17").Value = 50
$ws.Range("N17").Value = 3000
$ws.Range("O17").Value = 3000
$ws.Range("P17").Value = 3000
$ws.Range("R17").Value = "Provincia de Diguillín"
$ws.Range("S17").Value = 1500

# Row 18
$ws.Range("D18").Value = 44931
$ws.Range("M18").Value = 100

# Row 19
$ws.Range("D19").Value = 44966
$ws.Range("L19").Value = "Segunda"
$ws.Range("M19").Value = 30
$ws.Range("N19").Value = 2500
$ws.Range("O19").Value = 2500
$ws.Range("P19").Value = 2500
$ws.Range("S19").Value = 1250

# Row 20
$ws.Range("D20").Value = 44942
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 60

# Row 21
$ws.Range("D21").Value = 44937
$ws.Range("M21").Value = 100
$ws.Range("O21").Value = 3000
$ws.Range("P21").Value = 2750
$ws.Range("R21").Value = "Provincia de Diguillín"
$ws.Range("S21").Value = 1375

# Row 22
$ws.Range("D22").Value = 44953
$ws.Range("M22").Value = 30
$ws.Range("N22").Value = 3000
$ws.Range("O22").Value = 3000
$ws.Range("P22").Value = 3000
$ws.Range("S22").Value = 1500

# Row 23
$ws.Range("D23").Value = 44972
$ws.Range("L23").Value = "Segunda"
$ws.Range("M23").Value = 30
$ws.Range("N23").Value = 2500
$ws.Range("O23").Value = 2500
$ws.Range("P23").Value = 2500
$ws.Range("Q23").Value = "`$/bandeja 2 kilos"
$ws.Range("S23").Value = 1250
$ws.Range("T23").Value = 2

# Row 24
$ws.Range("D24").Value = 44952
$ws.Range("L24").Value = "Primera"
$ws.Range("N24").Value = 3000
$ws.Range("O24").Value = 3000
$ws.Range("P24").Value = 3000
$ws.Range("S24").Value = 1500

# Row 25
$ws.Range("D25").Value = 44596
$ws.Range("M25").Value = 120
$ws.Range("N25").Value = 2500
$ws.Range("O25").Value = 2700
$ws.Range("P25").Value = 2600
$ws.Range("R25").Value = "Provincia de Linares"
$ws.Range("S25").Value = 1300

# Row 26
$ws.Range("D26").Value = 44963
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 50
$ws.Range("N26").Value = 3000
$ws.Range("O26").Value = 3000
$ws.Range("P26").Value = 3000
$ws.Range("S26").Value = 1500

# Row 27
$ws.Range("D27").Value = 44963
$ws.Range("L27").Value = "Segunda"
$ws.Range("M27").Value = 50
$ws.Range("N27").Value = 2500
$ws.Range("O27").Value = 2500
$ws.Range("P27").Value = 2500
$ws.Range("S27").Value = 1250

# Row 28
$ws.Range("D28").Value = 44970
$ws.Range("M28").Value = 50
$ws.Range("N28").Value = 3000
$ws.Range("P28").Value = 3000
$ws.Range("S28").Value = 1500

# Row 29
$ws.Range("D29").Value = 44970
$ws.Range("L29").Value = "Segunda"
$ws.Range("M29").Value = 30
$ws.Range("N29").Value = 2500
$ws.Range("O29").Value = 2500
$ws.Range("P29").Value = 2500
$ws.Range("S29").Value = 1250

# Row 31
$ws.Range("D31").Value = 44960
$ws.Range("L31").Value = "Segunda"

# Row 32
$ws.Range("D32").Value = 44174
$ws.Range("M32").Value = 150
$ws.Range("N32").Value = 3700
$ws.Range("O32").Value = 3800
$ws.Range("P32").Value = 3747
$ws.Range("R32").Value = "Provincia de Linares"
$ws.Range("S32").Value = 1874

# --- Append new rows 33 and 34 ---
# Row 33
$ws.Range("A33").Value = 7
$ws.Range("B33").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C33").Value = "Ñuble"
$ws.Range("D33").Value = 44944
$ws.Range("E33").Value = 16
$ws.Range("F33").Value = "Fruta"
$ws.Range("G33").Value = 100101
$ws.Range("H33").Value = "Berries"
$ws.Range("I33").Value = 100101001
$ws.Range("J33").Value = "Arándano (blue)"
$ws.Range("K33").Value = "Sin especificar"
$ws.Range("L33").Value = "Primera"
$ws.Range("M33").Value = 60
$ws.Range("N33").Value = 2500
$ws.Range("O33").Value = 2500
$ws.Range("P33").Value = 2500
$ws.Range("Q33").Value = "`$/bandeja 2 kilos"
$ws.Range("R33").Value = "Provincia de Diguillín"
$ws.Range("S33").Value = 1250
$ws.Range("T33").Value = 2
$ws.Range("D33").NumberFormat = $ws.Range("D2").NumberFormat

# Row 34
$ws.Range("A34").Value = 7
$ws.Range("B34").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C34").Value = "Ñuble"
$ws.Range("D34").Value = 44539
$ws.Range("E34").Value = 16
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100101
$ws.Range("H34").Value = "Berries"
$ws.Range("I34").Value = 100101001
$ws.Range("J34").Value = "Arándano (blue)"
$ws.Range("K34").Value = "Sin especificar"
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 200
$ws.Range("N34").Value = 3800
$ws.Range("O34").Value = 4000
$ws.Range("P34").Value = 3900
$ws.Range("Q34").Value = "`$/bandeja 2 kilos"
$ws.Range("R34").Value = "Región del Maule"
$ws.Range("S34").Value = 1950
$ws.Range("T34").Value = 2
$ws.Range("D34").NumberFormat = $ws.Range("D2").NumberFormat
